$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Accredited)
$ws.Range("B2").Value = 15.346
$ws.Range("C2").Value = 23.623
$ws.Range("D2").Value = 1.726
$ws.Range("F2").Value = 40.695

# Row 3 (Unaccredited)
$ws.Range("B3").Value = 8.702
$ws.Range("C3").Value = 47.127
$ws.Range("D3").Value = 0.875
$ws.Range("E3").Value = 2.601
$ws.Range("F3").Value = 59.305

# Row 4 (COL_TOT)
$ws.Range("B4").Value = 24.048
$ws.Range("C4").Value = 70.75
$ws.Range("D4").Value = 2.601
$ws.Range("E4").Value = 2.601
